# Statistiques_M2_2010.xlsx - "amélioration de l'affichage des statistiques"
#
# Update the "Contenu du stage" pie-chart source data (rows 16-23, columns
# D/E/G on the Worksheet sheet): the student-count column (E) and the
# matching percentage labels (G) are refreshed with the real figures.
# The category labels in column D (C#, COBOL, C++, ASSEMBLEUR, ANDROID,
# JEE, DELPHI, PHP5) and the "Type entreprise" block (rows 25-28) keep
# their existing text - only their counts/percentages move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Percent-looking strings ("32 %", "0 %", ...) get auto-coerced to a
    # numeric percentage by the normal Value setter, so force the cell to
    # text first, write it, then drop the temporary number format again so
    # no stray style sticks around.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 16 - C#
$ws.Range("E16").Value = 8
Set-TextValue $ws.Range("G16") "32 %"

# Row 17 - COBOL
$ws.Range("E17").Value = 14
Set-TextValue $ws.Range("G17") "56 %"

# Row 18 - C++ (unchanged counts, kept explicit for clarity)
$ws.Range("E18").Value = 0
Set-TextValue $ws.Range("G18") "0 %"

# Row 19 - ASSEMBLEUR
$ws.Range("E19").Value = 2
Set-TextValue $ws.Range("G19") "8 %"

# Row 20 - ANDROID
$ws.Range("E20").Value = 1
Set-TextValue $ws.Range("G20") "4 %"

# Row 21 - JEE (unchanged)
$ws.Range("E21").Value = 0
Set-TextValue $ws.Range("G21") "0 %"

# Row 22 - DELPHI (unchanged)
$ws.Range("E22").Value = 0
Set-TextValue $ws.Range("G22") "0 %"

# Row 23 - PHP5 (unchanged)
$ws.Range("E23").Value = 0
Set-TextValue $ws.Range("G23") "0 %"

# "Type entreprise" block (rows 25-28) - values unchanged, left as-is.
